$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 39

# Columns A, B and D hold text that looks like a date / time / zero-padded
# number ("2024-01-09", "19:55:15", "01"). Excel would otherwise silently
# reinterpret these as a date serial, time serial or integer. Temporarily
# force a text number format while assigning the value, then clear the
# explicit formatting again so the new row ends up styled just like the
# other (unstyled) data rows.
$a = $ws.Cells.Item($row, 1)
$a.NumberFormat = "@"
$a.Value = "2024-01-09"
$a.ClearFormats()

$b = $ws.Cells.Item($row, 2)
$b.NumberFormat = "@"
$b.Value = "19:55:15"
$b.ClearFormats()

$ws.Cells.Item($row, 3).Value = "Tuesday"

$d = $ws.Cells.Item($row, 4)
$d.NumberFormat = "@"
$d.Value = "01"
$d.ClearFormats()

$ws.Cells.Item($row, 5).Value = 139483
$ws.Cells.Item($row, 6).Value = 142680
$ws.Cells.Item($row, 7).Value = 172111
$ws.Cells.Item($row, 8).Value = 147627
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 118468
$ws.Cells.Item($row, 11).Value = 224796
$ws.Cells.Item($row, 12).Value = 250773
$ws.Cells.Item($row, 13).Value = 185127
$ws.Cells.Item($row, 14).Value = 110385
$ws.Cells.Item($row, 15).Value = 40722
$ws.Cells.Item($row, 16).Value = 30866
$ws.Cells.Item($row, 17).Value = 72612
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42144
$ws.Cells.Item($row, 20).Value = -1
